$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

# Fill in the previously-empty row 7 with new ability data (no shifting of other rows)
$ws.Cells.Item(7, 1).Value = 50003
$ws.Cells.Item(7, 2).Value = "WEAPON_STING"
$ws.Cells.Item(7, 3).Value = "Weapon"
$ws.Cells.Item(7, 4).Value = "BE_STING"

$ws.Range("D13").Select()
